$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date values for rows 2-8 from 45185 to 45204
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45204
}
